$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8) beneath the existing prefixes table.
$ws.Cells.Item(8, 1).Value = "Class"
$ws.Cells.Item(8, 2).Value = $true
